$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 - header labels for the new "960 grid" table
$ws.Range("A13").Value = "Total"
$ws.Range("B13").Value = "Content"
$ws.Range("C13").Value = "Padding"
$ws.Range("D13").Value = "Margin"
$ws.Range("E13").Value = "Remainder"
$ws.Range("F13").Value = "Per"

# Row 14 - values / formulas for the grid calculations
$ws.Range("A14").Value = 960
$ws.Range("B14").Formula = "=142*5"
$ws.Range("C14").Value = 0
$ws.Range("D14").Formula = "=20*5"
$ws.Range("E14").Formula = "=A14-B14-C14-D14"
$ws.Range("F14").Formula = "=E14/5"

[void]$ws.Range("F14").Select()
